$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 17
$ws_ALC.Range("H17").Value = 844141.25
$ws_ALC.Range("J17").Value = 2059907.2
$ws_ALC.Range("L17").Value = 6179721.6
$ws_ALC.Range("N17").Value = -6180057.6

# ALC row 40
$ws_ALC.Range("H40").Value = 2776.182
$ws_ALC.Range("J40").Value = 2898.4
$ws_ALC.Range("L40").Value = 2898.4
$ws_ALC.Range("N40").Value = -3248.4

# ALC row 112
$ws_ALC.Range("H112").Value = 46173.348
$ws_ALC.Range("I112").Value = 0
$ws_ALC.Range("K112").Value = 0
$ws_ALC.Range("M112").ClearContents()

# ALC row 116
$ws_ALC.Range("H116").Value = 21242.9
$ws_ALC.Range("I116").Value = 19769.889
$ws_ALC.Range("K116").Value = 19769.889
$ws_ALC.Range("M116").Value = -16327.889

# ALC row 132
$ws_ALC.Range("H132").Value = 3126.5293
$ws_ALC.Range("I132").Value = 2321.6924
$ws_ALC.Range("K132").Value = 6965.0772
$ws_ALC.Range("M132").Value = -4435.0772

# ALC row 137
$ws_ALC.Range("H137").Value = 2356.5715
$ws_ALC.Range("I137").Value = 3099
$ws_ALC.Range("K137").Value = 9297
$ws_ALC.Range("M137").Value = -6747

# ALC row 138
$ws_ALC.Range("H138").Value = 13892757
$ws_ALC.Range("I138").Value = 1332.6666
$ws_ALC.Range("J138").Value = 20838468
$ws_ALC.Range("K138").Value = 3997.9998
$ws_ALC.Range("L138").Value = 62515404
$ws_ALC.Range("M138").Value = 1142.0002
$ws_ALC.Range("N138").Value = -62525684

# ARM row 88
$ws_ARM.Range("H88").Value = 2881.1428
$ws_ARM.Range("I88").Value = 2000
$ws_ARM.Range("J88").Value = 3028
$ws_ARM.Range("K88").Value = 2000
$ws_ARM.Range("L88").Value = 3028
$ws_ARM.Range("M88").Value = -1594
$ws_ARM.Range("N88").Value = -3840

# ARM row 91
$ws_ARM.Range("H91").Value = 2881.1428
$ws_ARM.Range("I91").Value = 2000
$ws_ARM.Range("J91").Value = 3028
$ws_ARM.Range("K91").Value = 2000
$ws_ARM.Range("L91").Value = 3028
$ws_ARM.Range("M91").Value = -596
$ws_ARM.Range("N91").Value = -5836

# ARM row 92
$ws_ARM.Range("H92").Value = 35000
$ws_ARM.Range("J92").Value = 35000
$ws_ARM.Range("L92").Value = 35000
$ws_ARM.Range("N92").Value = -39992

# ARM row 132
$ws_ARM.Range("H132").Value = 7618.278
$ws_ARM.Range("I132").Value = 7857.769
$ws_ARM.Range("J132").Value = 6995.6
$ws_ARM.Range("K132").Value = 23573.307
$ws_ARM.Range("L132").Value = 20986.8
$ws_ARM.Range("M132").Value = -21043.307
$ws_ARM.Range("N132").Value = -26046.8

# BSM row 134
$ws_BSM.Range("H134").Value = 2026.3334
$ws_BSM.Range("J134").Value = 1999.75
$ws_BSM.Range("L134").Value = 5999.25
$ws_BSM.Range("N134").Value = -11069.25

# CRP row 22
$ws_CRP.Range("H22").Value = 137.33333
$ws_CRP.Range("I22").Value = 137.33333
$ws_CRP.Range("K22").Value = 137.33333
$ws_CRP.Range("M22").Value = 212.66667

# CRP row 31
$ws_CRP.Range("H31").Value = 128458.625
$ws_CRP.Range("I31").Value = 253125
$ws_CRP.Range("K31").Value = 253125
$ws_CRP.Range("M31").Value = -252830

# CRP row 34
$ws_CRP.Range("H34").Value = 128458.625
$ws_CRP.Range("I34").Value = 253125
$ws_CRP.Range("K34").Value = 253125
$ws_CRP.Range("M34").Value = -252923

# CRP row 58
$ws_CRP.Range("H58").Value = 3056.7144
$ws_CRP.Range("I58").Value = 3080
$ws_CRP.Range("J58").Value = 2998.5
$ws_CRP.Range("K58").Value = 3080
$ws_CRP.Range("L58").Value = 2998.5
$ws_CRP.Range("M58").Value = -2877
$ws_CRP.Range("N58").Value = -3404.5

# CRP row 60
$ws_CRP.Range("H60").Value = 13242.25
$ws_CRP.Range("J60").Value = 18500
$ws_CRP.Range("L60").Value = 18500
$ws_CRP.Range("N60").Value = -19522

# CRP row 63
$ws_CRP.Range("H63").Value = 50001
$ws_CRP.Range("J63").Value = 50001
$ws_CRP.Range("L63").Value = 50001
$ws_CRP.Range("N63").Value = -51373

# CRP row 66
$ws_CRP.Range("H66").Value = 50001
$ws_CRP.Range("J66").Value = 50001
$ws_CRP.Range("L66").Value = 150003
$ws_CRP.Range("N66").Value = -156867

# CRP row 92
$ws_CRP.Range("H92").Value = 10000
$ws_CRP.Range("J92").Value = 10000
$ws_CRP.Range("L92").Value = 10000
$ws_CRP.Range("N92").Value = -14992

# CRP row 99
$ws_CRP.Range("H99").Value = 4931.65
$ws_CRP.Range("I99").Value = 4858.375
$ws_CRP.Range("K99").Value = 4858.375
$ws_CRP.Range("M99").Value = -3360.375

# CRP row 126
$ws_CRP.Range("H126").Value = 4931.65
$ws_CRP.Range("I126").Value = 4858.375
$ws_CRP.Range("K126").Value = 14575.125
$ws_CRP.Range("M126").Value = -12105.125

# CRP row 136
$ws_CRP.Range("H136").Value = 3056.7144
$ws_CRP.Range("I136").Value = 3080
$ws_CRP.Range("J136").Value = 2998.5
$ws_CRP.Range("K136").Value = 9240
$ws_CRP.Range("L136").Value = 8995.5
$ws_CRP.Range("M136").Value = -6690
$ws_CRP.Range("N136").Value = -14095.5

# CUL row 26
$ws_CUL.Range("H26").Value = 1205.4615
$ws_CUL.Range("I26").Value = 430.33334
$ws_CUL.Range("J26").Value = 1869.8572
$ws_CUL.Range("K26").Value = 1291.00002
$ws_CUL.Range("L26").Value = 5609.571599999999
$ws_CUL.Range("M26").Value = -1003.00002
$ws_CUL.Range("N26").Value = -6185.571599999999

# GSM row 92
$ws_GSM.Range("H92").Value = 9330.333000000001
$ws_GSM.Range("J92").Value = 8995.5
$ws_GSM.Range("L92").Value = 8995.5
$ws_GSM.Range("N92").Value = -12739.5

# GSM row 126
$ws_GSM.Range("H126").Value = 22543.166
$ws_GSM.Range("I126").Value = 26201.8
$ws_GSM.Range("J126").Value = 4250
$ws_GSM.Range("K126").Value = 78605.39999999999
$ws_GSM.Range("L126").Value = 12750
$ws_GSM.Range("M126").Value = -76135.39999999999
$ws_GSM.Range("N126").Value = -17690

# GSM row 132
$ws_GSM.Range("H132").Value = 4673.9165
$ws_GSM.Range("I132").Value = 3536.5
$ws_GSM.Range("J132").Value = 6948.75
$ws_GSM.Range("K132").Value = 10609.5
$ws_GSM.Range("L132").Value = 20846.25
$ws_GSM.Range("M132").Value = -8079.5
$ws_GSM.Range("N132").Value = -25906.25

# GSM row 136
$ws_GSM.Range("H136").Value = 31019.066
$ws_GSM.Range("J136").Value = 31019.066
$ws_GSM.Range("L136").Value = 93057.198
$ws_GSM.Range("N136").Value = -98157.198

# LTW row 7
$ws_LTW.Range("H7").Value = 5661.6665
$ws_LTW.Range("I7").Value = 5138.909
$ws_LTW.Range("K7").Value = 5138.909
$ws_LTW.Range("M7").Value = -5026.909

# LTW row 126
$ws_LTW.Range("H126").Value = 5661.6665
$ws_LTW.Range("I126").Value = 5138.909
$ws_LTW.Range("K126").Value = 15416.727
$ws_LTW.Range("M126").Value = -12946.727

# LTW row 132
$ws_LTW.Range("H132").Value = 5239.05
$ws_LTW.Range("I132").Value = 3941
$ws_LTW.Range("J132").Value = 7186.125
$ws_LTW.Range("K132").Value = 11823
$ws_LTW.Range("L132").Value = 21558.375
$ws_LTW.Range("M132").Value = -9293
$ws_LTW.Range("N132").Value = -26618.375

# LTW row 136
$ws_LTW.Range("H136").Value = 4816.952
$ws_LTW.Range("I136").Value = 4654
$ws_LTW.Range("J136").Value = 5142.857
$ws_LTW.Range("K136").Value = 13962
$ws_LTW.Range("L136").Value = 15428.571
$ws_LTW.Range("M136").Value = -11412
$ws_LTW.Range("N136").Value = -20528.571

# WVR row 122
$ws_WVR.Range("H122").Value = 2383.9524
$ws_WVR.Range("I122").Value = 2114.1765
$ws_WVR.Range("K122").Value = 6342.529500000001
$ws_WVR.Range("M122").Value = -3892.529500000001

# WVR row 132
$ws_WVR.Range("H132").Value = 1960.6666
$ws_WVR.Range("I132").Value = 1960.6666
$ws_WVR.Range("K132").Value = 5881.9998
$ws_WVR.Range("M132").Value = -3351.9998

# WVR row 136
$ws_WVR.Range("H136").Value = 2959.0833
$ws_WVR.Range("I136").Value = 2850.65
$ws_WVR.Range("J136").Value = 3501.25
$ws_WVR.Range("K136").Value = 8551.950000000001
$ws_WVR.Range("L136").Value = 10503.75
$ws_WVR.Range("M136").Value = -6001.950000000001
$ws_WVR.Range("N136").Value = -15603.75
